$d = $word.ActiveDocument

$oldSnippet = "Perioadele campaniei"
$newText = "Perioadele campaniei din Bootes: 14-23 mai, 13-22 iunie, 12-21 iulie"

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs($i)
    $full = $para.Range
    $text = $full.Text

    if ($text -like "*$oldSnippet*") {
        $start = $full.Start
        $end = $full.End - 1

        $target = $d.Range($start, $end)
        $target.Delete()

        $ins = $d.Range($start, $start)
        $ins.InsertAfter($newText)
    }
}
